# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# change: reorders / replaces the worker account-statement rows (B16:J33) with the
# updated dataset, keeping the same table shape, styles and formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
# H, I, J (Novedad Ingreso / Retiro / Observaciones) stay blank, as before.
$data = @(
    @("CC", "73182666",   "FABIAN PITALUA ZARZA",               "1802", 72000,  1142000),
    @("CC", "73182666",   "FABIAN PITALUA ZARZA",               "1801", 72000,  1142000),
    @("CC", "1143343026", "ANDRES FELIPE VASQUEZ MEJIA",        "1802", 96000,  2400000),
    @("CC", "1143343026", "ANDRES FELIPE VASQUEZ MEJIA",        "1801", 96000,  2400000),
    @("CC", "33069585",   "MARIA CLARA URIBE AGUILAR",          "1802", 73771,  1844292),
    @("CC", "33069585",   "MARIA CLARA URIBE AGUILAR",          "1801", 73771,  1844292),
    @("CC", "1143374517", "DARWIN CARIAGA GARCIA",              "1802", 36000,  900000),
    @("CC", "1143374517", "DARWIN CARIAGA GARCIA",              "1801", 36000,  900000),
    @("CC", "1143349287", "YENIFER PAOLA TATAR RODRIGUEZ",      "1802", 14667,  1000000),
    @("CC", "1143349287", "YENIFER PAOLA TATAR RODRIGUEZ",      "1801", 40000,  1000000),
    @("CE", "362441",     "LUIS FRANCISCO SAGARZAZU RODRIGUEZ", "1802", 29509,  737717),
    @("CE", "362441",     "LUIS FRANCISCO SAGARZAZU RODRIGUEZ", "1801", 29509,  737717),
    @("CE", "501276",     "MARCOS JOSE BORGES RAMOS",           "1802", 96000,  2400000),
    @("CE", "501276",     "MARCOS JOSE BORGES RAMOS",           "1801", 96000,  2400000),
    @("CC", "16787235",   "HECTOR FABIO FIGUEROA SOJET",        "1802", 240000, 6000000),
    @("CC", "16787235",   "HECTOR FABIO FIGUEROA SOJET",        "1801", 240000, 6000000),
    @("CC", "1127618941", "LEONARDO JOSE ROJAS LARA",           "1802", 96000,  2400000),
    @("CC", "1127618941", "LEONARDO JOSE ROJAS LARA",           "1801", 96000,  2400000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $ws.Cells.Item($row, 6).Value = $entry[4]
    $ws.Cells.Item($row, 7).Value = $entry[5]
}

# Columns B, C, E, G, H, I, J are "best fit" — let Excel recompute the
# optimal widths now that the underlying text/numbers have changed.
$ws.Range("B:C").EntireColumn.AutoFit() | Out-Null
$ws.Range("E:E").EntireColumn.AutoFit() | Out-Null
$ws.Range("G:J").EntireColumn.AutoFit() | Out-Null

$wb.Save()
